# Adds three more "logo" slides to the deck, derived from the existing
# title slide (slide 1): a dark logo-only slide, a light full (logo+text)
# slide, and a light logo-only slide.

$p = $ppt.ActivePresentation

# Color constants (PowerPoint RGB() packs as R + G*256 + B*65536).
$LIGHT_BG  = 16053235  # F3F3F4
$DARK_LINE = 4668723   # 333D47

# ---------------------------------------------------------------------
# Slide 2: duplicate of slide 1 (dark background, full logo + wordmark
# for now) - appended right after slide 1.
# ---------------------------------------------------------------------
$p.Slides.Item(1).Duplicate() | Out-Null
$s2 = $p.Slides.Item(2)

# ---------------------------------------------------------------------
# Slide 3: duplicate of slide 2 (still identical to slide 1 at this
# point), recolored to the light scheme: light background, dark
# "SIMPLY" text, dark logo outline (the "FI" text stays accent blue).
# ---------------------------------------------------------------------
$p.Slides.Item($p.Slides.Count).Duplicate() | Out-Null
$s3 = $p.Slides.Item(3)

$s3.Shapes.Item("Rechteck 4").Fill.ForeColor.RGB = $LIGHT_BG

$grp7_s3 = $s3.Shapes.Item("Gruppierung 7")
$tb_s3 = $grp7_s3.GroupItems.Item("Textfeld 21")
$simply_s3 = $tb_s3.TextFrame.TextRange.Characters(1, 6)
$simply_s3.Font.Color.RGB = $DARK_LINE

$grp1_s3 = $grp7_s3.GroupItems.Item("Gruppierung 1")
$grp1_s3.GroupItems.Item("Gleichschenkliges Dreieck 16").Line.ForeColor.RGB = $DARK_LINE
$grp1_s3.GroupItems.Item("Parallelogramm 8").Line.ForeColor.RGB = $DARK_LINE
$grp1_s3.GroupItems.Item("Gleichschenkliges Dreieck 23").Line.ForeColor.RGB = $DARK_LINE
$grp1_s3.GroupItems.Item("Gleichschenkliges Dreieck 24").Line.ForeColor.RGB = $DARK_LINE

# ---------------------------------------------------------------------
# Slide 4: duplicate of the now-recolored slide 3 (light background +
# dark logo), then strip the wordmark down to the logo only.
# ---------------------------------------------------------------------
$p.Slides.Item($p.Slides.Count).Duplicate() | Out-Null
$s4 = $p.Slides.Item(4)

$grp7_s4 = $s4.Shapes.Item("Gruppierung 7")
$grp7_s4.Ungroup() | Out-Null
$s4.Shapes.Item("Textfeld 21").Delete()

# ---------------------------------------------------------------------
# Now strip slide 2's wordmark too, leaving the dark logo only.
# ---------------------------------------------------------------------
$grp7_s2 = $s2.Shapes.Item("Gruppierung 7")
$grp7_s2.Ungroup() | Out-Null
$s2.Shapes.Item("Textfeld 21").Delete()

Write-Output ("Final slide count: " + $p.Slides.Count)
